$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.017.15'
$ws.Range("E2").Value = '  -0.45%  '

# Row 3
$ws.Range("D3").Value = '1.642.60'
$ws.Range("E3").Value = '  +0.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.68%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.21'
$ws.Range("E5").Value = '  -0.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("E6").Value = '  +0.83%  '

# Row 7
$ws.Range("E7").Value = '  +0.63%  '

# Row 8
$ws.Range("E8").Value = '  -0.24%  '

# Row 9
$ws.Range("E9").Value = '  +0.54%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  -0.42%  '

# Row 11
$ws.Range("E11").Value = '  +0.63%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.29'
$ws.Range("E12").Value = '  +0.11%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.869.16'
$ws.Range("E13").Value = '  -0.02%  '

# Row 14
$ws.Range("D14").Value = '1.629.55'
$ws.Range("E14").Value = '  -0.63%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.546'
$ws.Range("E15").Value = '  +0.23%  '

# Row 16
$ws.Range("E16").Value = '  +0.48%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.01'
$ws.Range("E17").Value = '  -0.35%  '

# Row 18
$ws.Range("D18").Value = '25.929.55'
$ws.Range("E18").Value = '  -0.80%  '

# Row 19
$ws.Range("E19").Value = '  +0.68%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.04'
$ws.Range("E20").Value = '  -0.82%  '

# Row 21
$ws.Range("E21").Value = '  -1.53%  '

# Row 22
$ws.Range("E22").Value = '  -0.90%  '

# Row 23
$ws.Range("E23").Value = '  -0.10%  '

# Row 24
$ws.Range("E24").Value = '  +1.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.132'
$ws.Range("E25").Value = '  +5.76%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.67'
$ws.Range("E26").Value = '  +1.52%  '

# Row 27
$ws.Range("E27").Value = '  +0.68%  '

# Row 28
$ws.Range("E28").Value = '  +0.05%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.57'
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("E30").Value = '  +0.03%  '

# Row 31
$ws.Range("E31").Value = '  -0.43%  '

# Row 32
$ws.Range("E32").Value = '  -1.88%  '

# Row 33
$ws.Range("E33").Value = '  +0.32%  '

# Row 35
$ws.Range("E35").Value = '  +2.56%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.904'
$ws.Range("E36").Value = '  -0.68%  '

# Row 37
$ws.Range("D37").Value = '1.134.68'
$ws.Range("E37").Value = '  +0.22%  '

# Row 38
$ws.Range("E38").Value = '  -1.37%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.47'
$ws.Range("E39").Value = '  -1.47%  '

# Row 40
$ws.Range("E40").Value = '  +0.34%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.52'
$ws.Range("E41").Value = '  +0.54%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.62'
$ws.Range("E42").Value = '  -0.68%  '

# Row 43
$ws.Range("E43").Value = '  +0.20%  '

# Row 44
$ws.Range("D44").Value = '1.778.37'
$ws.Range("E44").Value = '  -0.03%  '

# Row 45
$ws.Range("E45").Value = '  +3.52%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.77'
$ws.Range("E46").Value = '  +0.17%  '

# Row 47
$ws.Range("E47").Value = '  +2.74%  '

# Row 48
$ws.Range("E48").Value = '  -1.49%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.73'
$ws.Range("E49").Value = '  +0.84%  '

# Row 50
$ws.Range("E50").Value = '  -0.21%  '

# Row 51
$ws.Range("E51").Value = '  -0.33%  '

# Reset number format back to default (Normal style) for cells forced to text
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
